# Weekly update: a new price record (2022-07-xx, row id 7 / Terminal
# Hortofrutícola Agro Chillán / Cebollín) is inserted at the top of the
# data table, pushing the existing rows 4-16 down to 5-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 4 (shifts rows 4:16 -> 5:17,
# carrying formatting - e.g. the date-format style on column D - down
# with them).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44764
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112037
$ws.Range("G4").Value = "Cebollín"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("N4").Value = "$/docena de atados"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 2833
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = "Hortaliza"
